$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order: I=belongsTo, J=hasTopConcept, K=theme (altLabel column removed)
$ws.Range("I1").Value = "belongsTo"
$ws.Range("J1").Value = "hasTopConcept"
$ws.Range("K1").Value = "theme"

# Row 2 (erfpacht) - all null, unchanged values just shifted
$ws.Range("I2").Value = "null"
$ws.Range("J2").Value = "null"
$ws.Range("K2").Value = "null"

# Row 3 (opstal) - all null, unchanged values just shifted
$ws.Range("I3").Value = "null"
$ws.Range("J3").Value = "null"
$ws.Range("K3").Value = "null"

# Row 4 (volle_eigendom) - altLabel/theme data removed entirely
$ws.Range("I4").Value = "null"
$ws.Range("J4").Value = "null"
$ws.Range("K4").Value = "null"

# Row 5 (vruchtgebruik) - all null, unchanged values just shifted
$ws.Range("I5").Value = "null"
$ws.Range("J5").Value = "null"
$ws.Range("K5").Value = "null"

# Row 6 (conceptscheme) - belongsTo, hasTopConcept, theme values
$ws.Range("I6").Value = "https://data.omgeving.vlaanderen.be/id/dataset/codelijst-zakelijkrecht"
$ws.Range("J6").Value = "https://data.omgeving.vlaanderen.be/id/concept/zakelijkrecht/erfpacht|https://data.omgeving.vlaanderen.be/id/concept/zakelijkrecht/opstal|https://data.omgeving.vlaanderen.be/id/concept/zakelijkrecht/volle_eigendom|https://data.omgeving.vlaanderen.be/id/concept/zakelijkrecht/vruchtgebruik"
$ws.Range("K6").Value = "http://www.eionet.europa.eu/gemet/theme/5"

# Remove the now-unused column L entirely (was hasTopConcept, data moved to J)
$ws.Range("L1:L6").Delete()
